$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the style of the other header cells (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the corresponding numeric value in H2
$ws.Range("H2").Value = 0
